# Applies the numeric updates to the leve-profit tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), columns H:N, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(26, 8).Value = 16440
$ws.Cells.Item(26, 9).Value = 1013
$ws.Cells.Item(26, 10).Value = 55007.5
$ws.Cells.Item(26, 11).Value = 1013
$ws.Cells.Item(26, 12).Value = 55007.5
$ws.Cells.Item(26, 13).Value = -669
$ws.Cells.Item(26, 14).Value = -55695.5
$ws.Cells.Item(31, 8).Value = 1152.375
$ws.Cells.Item(31, 10).Value = 1170
$ws.Cells.Item(31, 12).Value = 3510
$ws.Cells.Item(31, 14).Value = -3970
$ws.Cells.Item(40, 8).Value = 1697.5
$ws.Cells.Item(40, 9).Value = 1050
$ws.Cells.Item(40, 10).Value = 2488.889
$ws.Cells.Item(40, 11).Value = 1050
$ws.Cells.Item(40, 12).Value = 2488.889
$ws.Cells.Item(40, 13).Value = -875
$ws.Cells.Item(40, 14).Value = -2838.889
$ws.Cells.Item(111, 8).Value = 1675.8
$ws.Cells.Item(111, 9).Value = 1859.6666
$ws.Cells.Item(111, 10).Value = 1400
$ws.Cells.Item(111, 11).Value = 5578.9998
$ws.Cells.Item(111, 12).Value = 4200
$ws.Cells.Item(111, 13).Value = -2511.9998
$ws.Cells.Item(111, 14).Value = -10334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 17859454
$ws.Cells.Item(2, 9).Value = 20834364
$ws.Cells.Item(2, 10).Value = 10000
$ws.Cells.Item(2, 11).Value = 20834364
$ws.Cells.Item(2, 12).Value = 10000
$ws.Cells.Item(2, 13).Value = -20834251
$ws.Cells.Item(2, 14).Value = -10226
$ws.Cells.Item(61, 8).Value = 1670.5
$ws.Cells.Item(61, 9).Value = 1190.96
$ws.Cells.Item(61, 10).Value = 5666.6665
$ws.Cells.Item(61, 11).Value = 1190.96
$ws.Cells.Item(61, 12).Value = 5666.6665
$ws.Cells.Item(61, 13).Value = -978.96
$ws.Cells.Item(61, 14).Value = -6090.6665
$ws.Cells.Item(116, 8).Value = 17859454
$ws.Cells.Item(116, 9).Value = 20834364
$ws.Cells.Item(116, 10).Value = 10000
$ws.Cells.Item(116, 11).Value = 20834364
$ws.Cells.Item(116, 12).Value = 10000
$ws.Cells.Item(116, 13).Value = -20832070
$ws.Cells.Item(116, 14).Value = -14588
$ws.Cells.Item(132, 8).Value = 38467536
$ws.Cells.Item(132, 9).Value = 66673224
$ws.Cells.Item(132, 11).Value = 200019672
$ws.Cells.Item(132, 13).Value = -200017142
$ws.Cells.Item(136, 8).Value = 1670.5
$ws.Cells.Item(136, 9).Value = 1190.96
$ws.Cells.Item(136, 10).Value = 5666.6665
$ws.Cells.Item(136, 11).Value = 3572.88
$ws.Cells.Item(136, 12).Value = 16999.9995
$ws.Cells.Item(136, 13).Value = -1022.88
$ws.Cells.Item(136, 14).Value = -22099.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 17859454
$ws.Cells.Item(3, 9).Value = 20834364
$ws.Cells.Item(3, 10).Value = 10000
$ws.Cells.Item(3, 11).Value = 20834364
$ws.Cells.Item(3, 12).Value = 10000
$ws.Cells.Item(3, 13).Value = -20834250
$ws.Cells.Item(3, 14).Value = -10228
$ws.Cells.Item(86, 8).Value = 1897.7059
$ws.Cells.Item(86, 9).Value = 1304.3572
$ws.Cells.Item(86, 11).Value = 1304.3572
$ws.Cells.Item(86, 13).Value = -181.3571999999999
$ws.Cells.Item(89, 8).Value = 1897.7059
$ws.Cells.Item(89, 9).Value = 1304.3572
$ws.Cells.Item(89, 11).Value = 6521.786
$ws.Cells.Item(89, 13).Value = -905.7860000000001
$ws.Cells.Item(99, 8).Value = 2167.1428
$ws.Cells.Item(99, 9).Value = 1448.75
$ws.Cells.Item(99, 10).Value = 3125
$ws.Cells.Item(99, 11).Value = 1448.75
$ws.Cells.Item(99, 12).Value = 3125
$ws.Cells.Item(99, 13).Value = 49.25
$ws.Cells.Item(99, 14).Value = -6121
$ws.Cells.Item(107, 8).Value = 1893.4117
$ws.Cells.Item(107, 9).Value = 1279.2
$ws.Cells.Item(107, 10).Value = 6500
$ws.Cells.Item(107, 11).Value = 1279.2
$ws.Cells.Item(107, 12).Value = 6500
$ws.Cells.Item(107, 13).Value = 640.8
$ws.Cells.Item(107, 14).Value = -10340

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2790
$ws.Cells.Item(16, 9).Value = 580
$ws.Cells.Item(16, 10).Value = 5000
$ws.Cells.Item(16, 11).Value = 580
$ws.Cells.Item(16, 12).Value = 5000
$ws.Cells.Item(16, 13).Value = -293
$ws.Cells.Item(16, 14).Value = -5574
$ws.Cells.Item(31, 8).Value = 3849709.2
$ws.Cells.Item(31, 9).Value = 5265997
$ws.Cells.Item(31, 11).Value = 5265997
$ws.Cells.Item(31, 13).Value = -5265702
$ws.Cells.Item(34, 8).Value = 3849709.2
$ws.Cells.Item(34, 9).Value = 5265997
$ws.Cells.Item(34, 11).Value = 5265997
$ws.Cells.Item(34, 13).Value = -5265795
$ws.Cells.Item(107, 8).Value = 2747.1538
$ws.Cells.Item(107, 9).Value = 1500
$ws.Cells.Item(107, 10).Value = 2851.0833
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 12).Value = 2851.0833
$ws.Cells.Item(107, 13).Value = 420
$ws.Cells.Item(107, 14).Value = -6691.0833
$ws.Cells.Item(113, 8).Value = 2790
$ws.Cells.Item(113, 9).Value = 580
$ws.Cells.Item(113, 10).Value = 5000
$ws.Cells.Item(113, 11).Value = 580
$ws.Cells.Item(113, 12).Value = 5000
$ws.Cells.Item(113, 13).Value = 1590
$ws.Cells.Item(113, 14).Value = -9340
$ws.Cells.Item(132, 8).Value = 3753.5173
$ws.Cells.Item(132, 9).Value = 2572.2104
$ws.Cells.Item(132, 11).Value = 7716.6312
$ws.Cells.Item(132, 13).Value = -5186.6312

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 788.7857
$ws.Cells.Item(107, 9).Value = 851.8182
$ws.Cells.Item(107, 10).Value = 748
$ws.Cells.Item(107, 11).Value = 2555.4546
$ws.Cells.Item(107, 12).Value = 2244
$ws.Cells.Item(107, 13).Value = -635.4546
$ws.Cells.Item(107, 14).Value = -6084
$ws.Cells.Item(109, 8).Value = 1445.1538
$ws.Cells.Item(109, 10).Value = 1719.7
$ws.Cells.Item(109, 12).Value = 5159.1
$ws.Cells.Item(109, 14).Value = -7239.1
$ws.Cells.Item(112, 8).Value = 3360
$ws.Cells.Item(115, 8).Value = 1127.25
$ws.Cells.Item(115, 9).Value = 623.6
$ws.Cells.Item(115, 10).Value = 1966.6666
$ws.Cells.Item(115, 11).Value = 1870.8
$ws.Cells.Item(115, 12).Value = 5899.9998
$ws.Cells.Item(115, 13).Value = -695.8000000000002
$ws.Cells.Item(115, 14).Value = -8249.9998
$ws.Cells.Item(118, 8).Value = 3870.5625
$ws.Cells.Item(118, 9).Value = 1929
$ws.Cells.Item(118, 10).Value = 4000
$ws.Cells.Item(118, 11).Value = 5787
$ws.Cells.Item(118, 12).Value = 12000
$ws.Cells.Item(118, 13).Value = -4544
$ws.Cells.Item(118, 14).Value = -14486
$ws.Cells.Item(121, 8).Value = 39250.125
$ws.Cells.Item(121, 9).Value = 433.33334
$ws.Cells.Item(121, 10).Value = 62540.2
$ws.Cells.Item(121, 11).Value = 1300.00002
$ws.Cells.Item(121, 12).Value = 187620.6
$ws.Cells.Item(121, 13).Value = 9.99998000000005
$ws.Cells.Item(121, 14).Value = -190240.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 4000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 10).Value = 10000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 12).Value = 10000
$ws.Cells.Item(107, 13).Value = 920
$ws.Cells.Item(107, 14).Value = -13840

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 771.5714
$ws.Cells.Item(16, 9).Value = 716.8333
$ws.Cells.Item(16, 10).Value = 1100
$ws.Cells.Item(16, 11).Value = 716.8333
$ws.Cells.Item(16, 12).Value = 1100
$ws.Cells.Item(16, 13).Value = -546.8333
$ws.Cells.Item(16, 14).Value = -1440

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 45000
$ws.Cells.Item(46, 10).Value = 45000
$ws.Cells.Item(46, 12).Value = 45000
$ws.Cells.Item(46, 14).Value = -45462
$ws.Cells.Item(81, 8).Value = 1500
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 1500
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 12).Value = 3000
$ws.Cells.Item(81, 13).ClearContents()
$ws.Cells.Item(81, 14).Value = -5122
$ws.Cells.Item(84, 8).Value = 1500
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 1500
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 15000
$ws.Cells.Item(84, 13).ClearContents()
$ws.Cells.Item(84, 14).Value = -25608
$ws.Cells.Item(107, 8).Value = 2457.5715
$ws.Cells.Item(107, 9).Value = 2036.909
$ws.Cells.Item(107, 10).Value = 4000
$ws.Cells.Item(107, 11).Value = 6110.727000000001
$ws.Cells.Item(107, 12).Value = 12000
$ws.Cells.Item(107, 13).Value = -4190.727000000001
$ws.Cells.Item(107, 14).Value = -15840
$ws.Cells.Item(132, 8).Value = 255388.67
$ws.Cells.Item(132, 9).Value = 335718.12
$ws.Cells.Item(132, 10).Value = 14400.3
$ws.Cells.Item(132, 11).Value = 1007154.36
$ws.Cells.Item(132, 12).Value = 43200.89999999999
$ws.Cells.Item(132, 13).Value = -1004624.36
$ws.Cells.Item(132, 14).Value = -48260.89999999999
$ws.Cells.Item(134, 8).Value = 45000
$ws.Cells.Item(134, 10).Value = 45000
$ws.Cells.Item(134, 12).Value = 135000
$ws.Cells.Item(134, 14).Value = -140070
